# Feria Lagunitas de Puerto Montt - Uva
# Insert a new weekly price record at row 246, pushing the existing
# rows (246-287) down by one (becoming 247-288).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 246 (shifts 246:287 -> 247:288)
$ws.Rows.Item(246).Insert()

# Populate the new row 246 with the new record's data
$ws.Cells.Item(246, 1).Value  = 4
$ws.Cells.Item(246, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(246, 3).Value  = "Los Lagos"
$ws.Cells.Item(246, 4).Value  = 44951
$ws.Cells.Item(246, 5).Value  = 10
$ws.Cells.Item(246, 6).Value  = "Fruta"
$ws.Cells.Item(246, 7).Value  = 100109
$ws.Cells.Item(246, 8).Value  = "Uva"
$ws.Cells.Item(246, 9).Value  = 100109001
$ws.Cells.Item(246, 10).Value = "Uva"
$ws.Cells.Item(246, 11).Value = "Red Globe"
$ws.Cells.Item(246, 12).Value = "Primera"
$ws.Cells.Item(246, 13).Value = 200
$ws.Cells.Item(246, 14).Value = 12000
$ws.Cells.Item(246, 15).Value = 13000
$ws.Cells.Item(246, 16).Value = 12500
$ws.Cells.Item(246, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(246, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(246, 19).Value = 1562
$ws.Cells.Item(246, 20).Value = 8

# Match the date cell number format used by the other rows of column D
# so the new row renders consistently with the rest of the table.
$ws.Cells.Item(246, 4).NumberFormat = $ws.Cells.Item(247, 4).NumberFormat()
